$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Ports sheet: bump charging-rate column (E) for existing ports, and add
# a new Port 6 row (row 7) mirroring the formatting of the row above it.
# ---------------------------------------------------------------------
$wsPorts = $wb.Worksheets.Item("Ports")

$wsPorts.Range("A6:F6").Copy()
$wsPorts.Range("A7:F7").PasteSpecial(-4122)

$wsPorts.Range("E2").Value = 5
$wsPorts.Range("E3").Value = 5
$wsPorts.Range("E4").Value = 5
$wsPorts.Range("E5").Value = 5
$wsPorts.Range("E6").Value = 5

$wsPorts.Range("A7").Value = 6
$wsPorts.Range("B7").Value = 50
$wsPorts.Range("C7").Value = 72
$wsPorts.Range("D7").Value = 0.6
$wsPorts.Range("E7").Value = 5
$wsPorts.Range("F7").Value = 20

# ---------------------------------------------------------------------
# Operator sheet: add a 6th aircraft / datalink slot (column H) to the
# "Serviced Ports" / "Charging Equipment" selector rows, widen the
# Landing-Slots COUNTIF for row 13 to include it, and bump the fleet
# size inputs (C10/D10).
# ---------------------------------------------------------------------
$wsOp = $wb.Worksheets.Item("Operator")

$wsOp.Range("C10").Value = 8
$wsOp.Range("D10").Value = 8

$wsOp.Range("B13").Formula = "=COUNTIF(C13:I13,""Yes"")"

$wsOp.Range("C13").Copy()
$wsOp.Range("H13").PasteSpecial(-4122)
$wsOp.Range("H13").Value = "Yes"

$wsOp.Range("C14").Copy()
$wsOp.Range("H14").PasteSpecial(-4122)
$wsOp.Range("H14").Value = "Slow"

$wsOp.Range("C14:H14").Validation.Delete()
$wsOp.Range("C14:H14").Validation.Add(3, 1, 1, """Slow, Fast, None""")

# Make "Operator" the active/visible tab (was "Aircraft").
$wsOp.Activate()
[void]$wsOp.Range("A1").Select()
